$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price values to preserve text type
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply cell value updates per diff
$ws.Range("D2").Value = "26.696.02"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.634.32"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "216.99"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "0.249"
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("D9").Value = "0.0621"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").Value = "18.97"
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "1.862.51"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").Value = "1.636.26"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").Value = "0.523"
$ws.Range("E15").Value = "  -1.81%  "
$ws.Range("D16").Value = "64.28"
$ws.Range("E16").Value = "  -1.76%  "
$ws.Range("D17").Value = "26.661.62"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "0.0₃0724"
$ws.Range("E18").Value = "  -2.67%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "210.13"
$ws.Range("E20").Value = "  -3.95%  "
$ws.Range("D21").Value = "4.32"
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "6.17"
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").Value = "2.32"
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("D24").Value = "9.22"
$ws.Range("E24").Value = "  -2.97%  "
$ws.Range("D25").Value = "145.51"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("E27").Value = "  -2.49%  "
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("D29").Value = "15.51"
$ws.Range("E29").Value = "  -1.43%  "
$ws.Range("D30").Value = "0.0502"
$ws.Range("E30").Value = "  -3.00%  "
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("D33").Value = "2.96"
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("D34").Value = "1.271.65"
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("E35").Value = "  -1.96%  "
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("E37").Value = "  -2.29%  "
$ws.Range("D38").Value = "0.527"
$ws.Range("E38").Value = "  -1.79%  "
$ws.Range("D39").Value = "0.805"
$ws.Range("E39").Value = "  -2.48%  "
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("E41").Value = "  -1.73%  "
$ws.Range("E42").Value = "  -2.39%  "
$ws.Range("D43").Value = "1.773.62"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("E44").Value = "  -3.74%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "91.03"
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "59.98"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  -2.72%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0102"
$ws.Range("E48").Value = "  -3.21%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0517"
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "7.52"
$ws.Range("E50").Value = "  -3.23%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.0958"
$ws.Range("E51").Value = "  -1.06%  "
